# 17 veebruar kodutoo exception parandus
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nädal 4")

# Fill in row 8 (second log entry: 19 Feb 2020, 08:00-10:00, 60 min, praktikum)
$ws.Range("B8").Value = 43880
$ws.Range("C8").Value = 0.33333333333333331
$ws.Range("D8").Value = 0.41666666666666669
$ws.Range("E8").Value = "-"
$ws.Range("F8").Value = 60
$ws.Range("G8").Value = "Praktikum breakpoint ja debuggimine"
$ws.Range("I8").Value = "x"

# Move the active selection to I8 to match the saved workbook view
$ws.Range("I8").Select()
